$p = $ppt.ActivePresentation

# The deck currently carries the "Integral" theme colours on the live
# slide-master theme part (reached through Slide.ThemeColorScheme).
# Reset that colour scheme back to the standard "Office Theme" palette
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink), matching the colours
# that the notes-master's (unused) theme already carried.

function RGBVal($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = RGBVal($officeColors[$i - 1])
}
